# Apply 2023-11-24 daily crime count update (violent-crime-full-year.xlsx)
# Updates column J (year 2023) totals across the Citywide, By Neighborhood,
# and per-neighborhood worksheets.

$wb = $excel.ActiveWorkbook

# 1: Citywide Totals
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 10).Value = 6911
$ws.Cells.Item(3, 10).Value = 7309
$ws.Cells.Item(4, 10).Value = 1589
$ws.Cells.Item(6, 10).Value = 9803
$ws.Cells.Item(7, 10).Value = 26187

# 2: By Neighborhood
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 10).Value = 208
$ws.Cells.Item(4, 10).Value = 116
$ws.Cells.Item(8, 10).Value = 1646
$ws.Cells.Item(10, 10).Value = 193
$ws.Cells.Item(11, 10).Value = 456
$ws.Cells.Item(15, 10).Value = 313
$ws.Cells.Item(18, 10).Value = 218
$ws.Cells.Item(19, 10).Value = 762
$ws.Cells.Item(20, 10).Value = 551
$ws.Cells.Item(23, 10).Value = 241
$ws.Cells.Item(29, 10).Value = 1413
$ws.Cells.Item(31, 10).Value = 266
$ws.Cells.Item(33, 10).Value = 1183
$ws.Cells.Item(35, 10).Value = 31
$ws.Cells.Item(37, 10).Value = 813
$ws.Cells.Item(41, 10).Value = 183
$ws.Cells.Item(42, 10).Value = 1130
$ws.Cells.Item(43, 10).Value = 225
$ws.Cells.Item(45, 10).Value = 38
$ws.Cells.Item(49, 10).Value = 164
$ws.Cells.Item(52, 10).Value = 670
$ws.Cells.Item(53, 10).Value = 382
$ws.Cells.Item(54, 10).Value = 505
$ws.Cells.Item(55, 10).Value = 411
$ws.Cells.Item(57, 10).Value = 119
$ws.Cells.Item(59, 10).Value = 29
$ws.Cells.Item(60, 10).Value = 153
$ws.Cells.Item(62, 10).Value = 9
$ws.Cells.Item(63, 10).Value = 83
$ws.Cells.Item(65, 10).Value = 658
$ws.Cells.Item(67, 10).Value = 976
$ws.Cells.Item(71, 10).Value = 85
$ws.Cells.Item(76, 10).Value = 376
$ws.Cells.Item(77, 10).Value = 182
$ws.Cells.Item(78, 10).Value = 307
$ws.Cells.Item(79, 10).Value = 734
$ws.Cells.Item(83, 10).Value = 524
$ws.Cells.Item(85, 10).Value = 1077
$ws.Cells.Item(86, 10).Value = 166
$ws.Cells.Item(87, 10).Value = 85
$ws.Cells.Item(88, 10).Value = 278
$ws.Cells.Item(89, 10).Value = 331
$ws.Cells.Item(91, 10).Value = 303
$ws.Cells.Item(93, 10).Value = 109
$ws.Cells.Item(94, 10).Value = 282
$ws.Cells.Item(95, 10).Value = 383
$ws.Cells.Item(98, 10).Value = 195
$ws.Cells.Item(99, 10).Value = 399
$ws.Cells.Item(101, 10).Value = 26187

# 6: Belmont Cragin
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(2, 10).Value = 130
$ws.Cells.Item(6, 10).Value = 211
$ws.Cells.Item(7, 10).Value = 456

# 7: Uptown
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(3, 10).Value = 96
$ws.Cells.Item(6, 10).Value = 98
$ws.Cells.Item(7, 10).Value = 331

# 8: South Shore
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(3, 10).Value = 385
$ws.Cells.Item(6, 10).Value = 310
$ws.Cells.Item(7, 10).Value = 1077

# 9: Little Village
$ws = $wb.Worksheets.Item(9)
$ws.Cells.Item(3, 10).Value = 190
$ws.Cells.Item(6, 10).Value = 288
$ws.Cells.Item(7, 10).Value = 670

# 11: Logan Square
$ws = $wb.Worksheets.Item(11)
$ws.Cells.Item(2, 10).Value = 67
$ws.Cells.Item(6, 10).Value = 253
$ws.Cells.Item(7, 10).Value = 382

# 12: Austin
$ws = $wb.Worksheets.Item(12)
$ws.Cells.Item(2, 10).Value = 440
$ws.Cells.Item(3, 10).Value = 490
$ws.Cells.Item(6, 10).Value = 586
$ws.Cells.Item(7, 10).Value = 1646

# 13: South Chicago
$ws = $wb.Worksheets.Item(13)
$ws.Cells.Item(2, 10).Value = 154
$ws.Cells.Item(7, 10).Value = 524

# 14: Garfield Park
$ws = $wb.Worksheets.Item(14)
$ws.Cells.Item(3, 10).Value = 393
$ws.Cells.Item(7, 10).Value = 1183

# 15: West Pullman
$ws = $wb.Worksheets.Item(15)
$ws.Cells.Item(2, 10).Value = 134
$ws.Cells.Item(4, 10).Value = 19
$ws.Cells.Item(7, 10).Value = 383

# 16: Grand Crossing
$ws = $wb.Worksheets.Item(16)
$ws.Cells.Item(3, 10).Value = 275
$ws.Cells.Item(7, 10).Value = 813

# 17: New City
$ws = $wb.Worksheets.Item(17)
$ws.Cells.Item(2, 10).Value = 189
$ws.Cells.Item(6, 10).Value = 242
$ws.Cells.Item(7, 10).Value = 658

# 18: Woodlawn
$ws = $wb.Worksheets.Item(18)
$ws.Cells.Item(5, 10).Value = 12
$ws.Cells.Item(6, 10).Value = 102
$ws.Cells.Item(7, 10).Value = 399

# 20: Gage Park
$ws = $wb.Worksheets.Item(20)
$ws.Cells.Item(6, 10).Value = 90
$ws.Cells.Item(7, 10).Value = 266

# 21: North Lawndale
$ws = $wb.Worksheets.Item(21)
$ws.Cells.Item(2, 10).Value = 248
$ws.Cells.Item(7, 10).Value = 976

# 23: Lincoln Park
$ws = $wb.Worksheets.Item(23)
$ws.Cells.Item(2, 10).Value = 30
$ws.Cells.Item(7, 10).Value = 164

# 24: Loop
$ws = $wb.Worksheets.Item(24)
$ws.Cells.Item(2, 10).Value = 125
$ws.Cells.Item(6, 10).Value = 236
$ws.Cells.Item(7, 10).Value = 505

# 25: Englewood
$ws = $wb.Worksheets.Item(25)
$ws.Cells.Item(3, 10).Value = 502
$ws.Cells.Item(6, 10).Value = 359
$ws.Cells.Item(7, 10).Value = 1413

# 27: Chatham
$ws = $wb.Worksheets.Item(27)
$ws.Cells.Item(2, 10).Value = 184
$ws.Cells.Item(6, 10).Value = 295
$ws.Cells.Item(7, 10).Value = 762

# 29: River North
$ws = $wb.Worksheets.Item(29)
$ws.Cells.Item(6, 10).Value = 202
$ws.Cells.Item(7, 10).Value = 376

# 31: Hermosa
$ws = $wb.Worksheets.Item(31)
$ws.Cells.Item(2, 10).Value = 36
$ws.Cells.Item(7, 10).Value = 183

# 32: Humboldt Park
$ws = $wb.Worksheets.Item(32)
$ws.Cells.Item(3, 10).Value = 222
$ws.Cells.Item(6, 10).Value = 602
$ws.Cells.Item(7, 10).Value = 1130

# 34: Avondale
$ws = $wb.Worksheets.Item(34)
$ws.Cells.Item(6, 10).Value = 109
$ws.Cells.Item(7, 10).Value = 193

# 35: Rogers Park
$ws = $wb.Worksheets.Item(35)
$ws.Cells.Item(3, 10).Value = 96
$ws.Cells.Item(7, 10).Value = 307

# 36: Lower West Side
$ws = $wb.Worksheets.Item(36)
$ws.Cells.Item(6, 10).Value = 231
$ws.Cells.Item(7, 10).Value = 411

# 39: Douglas
$ws = $wb.Worksheets.Item(39)
$ws.Cells.Item(2, 10).Value = 65
$ws.Cells.Item(7, 10).Value = 241

# 40: Washington Park
$ws = $wb.Worksheets.Item(40)
$ws.Cells.Item(6, 10).Value = 77
$ws.Cells.Item(7, 10).Value = 303

# 42: Roseland
$ws = $wb.Worksheets.Item(42)
$ws.Cells.Item(2, 10).Value = 204
$ws.Cells.Item(3, 10).Value = 247
$ws.Cells.Item(7, 10).Value = 734

# 44: Chicago Lawn
$ws = $wb.Worksheets.Item(44)
$ws.Cells.Item(2, 10).Value = 154
$ws.Cells.Item(6, 10).Value = 156
$ws.Cells.Item(7, 10).Value = 551

# 45: Calumet Heights
$ws = $wb.Worksheets.Item(45)
$ws.Cells.Item(6, 10).Value = 101
$ws.Cells.Item(7, 10).Value = 218

# 48: West Lawn
$ws = $wb.Worksheets.Item(48)
$ws.Cells.Item(2, 10).Value = 29
$ws.Cells.Item(7, 10).Value = 109

# 50: Garfield Ridge
$ws = $wb.Worksheets.Item(50)
$ws.Cells.Item(4, 10).Value = 6
$ws.Cells.Item(6, 10).Value = 45

# 51: West Loop
$ws = $wb.Worksheets.Item(51)
$ws.Cells.Item(3, 10).Value = 55
$ws.Cells.Item(7, 10).Value = 282

# 54: Brighton Park
$ws = $wb.Worksheets.Item(54)
$ws.Cells.Item(2, 10).Value = 90
$ws.Cells.Item(6, 10).Value = 136
$ws.Cells.Item(7, 10).Value = 313

# 55: Wicker Park
$ws = $wb.Worksheets.Item(55)
$ws.Cells.Item(6, 10).Value = 126
$ws.Cells.Item(7, 10).Value = 195

# 60: Gold Coast
$ws = $wb.Worksheets.Item(60)
$ws.Cells.Item(3, 10).Value = 5
$ws.Cells.Item(7, 10).Value = 31

# 63: Montclare
$ws = $wb.Worksheets.Item(63)
$ws.Cells.Item(3, 10).Value = 6
$ws.Cells.Item(7, 10).Value = 29

# 64: Albany Park
$ws = $wb.Worksheets.Item(64)
$ws.Cells.Item(2, 10).Value = 61
$ws.Cells.Item(6, 10).Value = 78
$ws.Cells.Item(7, 10).Value = 208

# 68: United Center
$ws = $wb.Worksheets.Item(68)
$ws.Cells.Item(3, 10).Value = 67
$ws.Cells.Item(6, 10).Value = 141
$ws.Cells.Item(7, 10).Value = 278

# 72: Streeterville
$ws = $wb.Worksheets.Item(72)
$ws.Cells.Item(4, 10).Value = 88
$ws.Cells.Item(7, 10).Value = 166

# 77: Mckinley Park
$ws = $wb.Worksheets.Item(77)
$ws.Cells.Item(6, 10).Value = 53
$ws.Cells.Item(7, 10).Value = 119

# 78: Morgan Park
$ws = $wb.Worksheets.Item(78)
$ws.Cells.Item(3, 10).Value = 43
$ws.Cells.Item(7, 10).Value = 153

# 79: Hyde Park
$ws = $wb.Worksheets.Item(79)
$ws.Cells.Item(4, 10).Value = 20
$ws.Cells.Item(7, 10).Value = 225

# 81: Oakland
$ws = $wb.Worksheets.Item(81)
$ws.Cells.Item(3, 10).Value = 24
$ws.Cells.Item(7, 10).Value = 85

# 84: Riverdale
$ws = $wb.Worksheets.Item(84)
$ws.Cells.Item(6, 10).Value = 33
$ws.Cells.Item(7, 10).Value = 182

# 85: Jackson Park
$ws = $wb.Worksheets.Item(85)
$ws.Cells.Item(2, 10).Value = 14
$ws.Cells.Item(7, 10).Value = 38

# 90: Archer Heights
$ws = $wb.Worksheets.Item(90)
$ws.Cells.Item(2, 10).Value = 37
$ws.Cells.Item(7, 10).Value = 116

# 92: Ukrainian Village
$ws = $wb.Worksheets.Item(92)
$ws.Cells.Item(6, 10).Value = 58
$ws.Cells.Item(7, 10).Value = 85

# 97: Millenium Park
$ws = $wb.Worksheets.Item(97)
$ws.Cells.Item(4, 10).Value = 2

# 98: Museum Campus
$ws = $wb.Worksheets.Item(98)
$ws.Cells.Item(7, 10).Value = 9

